$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header fields ---
$ws.Range("C2").Value = 45364

# The "Ficha" number is a text identifier, not a quantity, so it must stay
# string-typed even though it looks numeric. A direct .Value assignment of a
# pure-digit string gets auto-coerced to a number by Excel, so instead enter
# it as a text formula and immediately collapse it to a static value with
# Copy + PasteSpecial(xlPasteValues) - this keeps the stored type as text
# without leaving a formula behind or touching any cell's number format.
$ws.Range("C3").Formula = "=""2556845"""
$ws.Range("C3").Copy()
$ws.Range("C3").PasteSpecial(-4163)
# The paste-special drops the pre-existing C3:F3 merge, so restore it.
$ws.Range("C3:F3").Merge()

$ws.Range("C4").Value = "DISEÑO E INTEGRACION DE MULTIMEDIA"
$ws.Range("C6").Value = 44760
$ws.Range("C7").Value = 45124

# --- Replace / extend the student listing starting at row 11 ---
$rows = @(
    @("CC", 1005178211, "JUAN CAMILO",     "DELGADO CARRASCAL", "CERTIFICADO"),
    @("CC", 1005181992, "JOHAN",           "VARGAS CALDERIN",   "POR CERTIFICAR"),
    @("CC", 1005184329, "WILLIAM ANDRES",  "LOPEZ RIOS",         "CANCELADO"),
    @("CC", 1005185919, "SEBASTIAN",       "PERTUZ SAMPAYO",     "CERTIFICADO"),
    @("CC", 1005220651, "BRAYAN EDUARDO",  "BADILLO HERRERA",    "CERTIFICADO"),
    @("CC", 1005239745, "SARAY DUVIANA",   "UNRIZA JAIMES",      "CERTIFICADO"),
    @("CC", 1005241421, "CLARA LUCIA",     "RUIZ MONSALVE",      "RETIRO VOLUNTARIO"),
    @("CC", 1043962939, "DANNA KAROLAY",   "RESTREPO SOSA",      "CERTIFICADO"),
    @("CC", 1048457729, "DAYANA",          "URRUCHURTU NIÑO",    "CERTIFICADO"),
    @("TI", 1049019898, "KAREN YURLEIDY",  "MARIN VARGAS",       "RETIRO VOLUNTARIO"),
    @("CC", 1087985197, "GISELL MARIANA",  "MARIN LARROTA",      "CERTIFICADO"),
    @("CC", 1096184002, "DANIELA",         "ROJAS BOTELLO",      "CERTIFICADO"),
    @("CC", 1096186262, "KEVIN ANDRES",    "PARADA SUAREZ",      "RETIRO VOLUNTARIO"),
    @("CC", 1096189477, "KAMILA",          "QUINTERO CARREÑO",   "CERTIFICADO"),
    @("CC", 1097183074, "MARIA JOSE",      "ORTIZ GUIZA",        "CERTIFICADO"),
    @("CC", 1144182405, "CAROLAIN",        "ABANIS PEREZ",       "CERTIFICADO"),
    @("CC", 63469380,   "VIDA EMPERATRIZ", "SANTOS YAIN",        "CERTIFICADO")
)

$r = 11
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
